# Updated cryptos list on Mon Feb 12 03:42:57 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'48.160.72"
$ws.Range("E2").Value = "'  -0.48%  "

# Row 3
$ws.Range("D3").Value = "'2.500.07"
$ws.Range("E3").Value = "'  -0.63%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.07%  "

# Row 5
$ws.Range("D5").Value = "'320.03"
$ws.Range("E5").Value = "'  -1.13%  "

# Row 6
$ws.Range("D6").Value = "'106.49"
$ws.Range("E6").Value = "'  -2.73%  "

# Row 7
$ws.Range("D7").Value = "'0.524"
$ws.Range("E7").Value = "'  -0.43%  "

# Row 8
$ws.Range("E8").Value = "'  +0.05%  "

# Row 9
$ws.Range("D9").Value = "'0.540"
$ws.Range("E9").Value = "'  -4.07%  "

# Row 10
$ws.Range("D10").Value = "'38.85"
$ws.Range("E10").Value = "'  -4.18%  "

# Row 11
$ws.Range("D11").Value = "'19.94"
$ws.Range("E11").Value = "'  +1.42%  "

# Row 12
$ws.Range("E12").Value = "'  -1.80%  "

# Row 13
$ws.Range("E13").Value = "'  -0.46%  "

# Row 14
$ws.Range("D14").Value = "'7.08"
$ws.Range("E14").Value = "'  -2.05%  "

# Row 15
$ws.Range("D15").Value = "'2.891.78"
$ws.Range("E15").Value = "'  -0.47%  "

# Row 16
$ws.Range("D16").Value = "'2.496.91"
$ws.Range("E16").Value = "'  -0.69%  "

# Row 17
$ws.Range("D17").Value = "'0.834"
$ws.Range("E17").Value = "'  -2.41%  "

# Row 18
$ws.Range("D18").Value = "'48.035.59"
$ws.Range("E18").Value = "'  -0.33%  "

# Row 19
$ws.Range("B19").Value = "'ImmutableX"
$ws.Range("C19").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").Value = "'2.99"
$ws.Range("E19").Value = "'  +8.86%  "

# Row 20
$ws.Range("B20").Value = "'InternetComputer(DFINITY)"
$ws.Range("C20").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'13.02"
$ws.Range("E20").Value = "'  -2.86%  "

# Row 21
$ws.Range("E21").Value = "'  +0.18%  "

# Row 22
$ws.Range("D22").Value = "'0.0₃0935"
$ws.Range("E22").Value = "'  -1.23%  "

# Row 23
$ws.Range("D23").Value = "'71.18"
$ws.Range("E23").Value = "'  -0.60%  "

# Row 24
$ws.Range("D24").Value = "'272.81"
$ws.Range("E24").Value = "'  +1.57%  "

# Row 25
$ws.Range("E25").Value = "'  -1.94%  "

# Row 26
$ws.Range("E26").Value = "'  -0.07%  "

# Row 27
$ws.Range("D27").Value = "'25.83"
$ws.Range("E27").Value = "'  -1.03%  "

# Row 28
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = "'  +2.88%  "

# Row 29
$ws.Range("D29").Value = "'0.142"
$ws.Range("E29").Value = "'  -0.45%  "

# Row 30
$ws.Range("D30").Value = "'9.74"
$ws.Range("E30").Value = "'  -4.67%  "

# Row 31
$ws.Range("D31").Value = "'35.04"
$ws.Range("E31").Value = "'  -0.92%  "

# Row 32
$ws.Range("E32").Value = "'  -0.96%  "

# Row 33
$ws.Range("B33").Value = "'FirstDigitalUSD"
$ws.Range("C33").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "'  -0.18%  "

# Row 34
$ws.Range("B34").Value = "'Celestia"
$ws.Range("C34").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "'19.13"
$ws.Range("E34").Value = "'  -4.75%  "

# Row 35
$ws.Range("E35").Value = "'  -2.09%  "

# Row 36
$ws.Range("D36").Value = "'0.0777"
$ws.Range("E36").Value = "'  -1.49%  "

# Row 37
$ws.Range("E37").Value = "'  -2.22%  "

# Row 38
$ws.Range("D38").Value = "'4.60"
$ws.Range("E38").Value = "'  -2.46%  "

# Row 39
$ws.Range("E39").Value = "'  -3.36%  "

# Row 40
$ws.Range("E40").Value = "'  -1.49%  "

# Row 41
$ws.Range("D41").Value = "'121.24"
$ws.Range("E41").Value = "'  +1.83%  "

# Row 42
$ws.Range("D42").Value = "'22.09"
$ws.Range("E42").Value = "'  -0.58%  "

# Row 43
$ws.Range("E43").Value = "'  +0.87%  "

# Row 44
$ws.Range("D44").Value = "'0.0305"
$ws.Range("E44").Value = "'  +1.94%  "

# Row 45
$ws.Range("D45").Value = "'2.010.96"
$ws.Range("E45").Value = "'  +0.44%  "

# Row 46
$ws.Range("E46").Value = "'  +2.10%  "

# Row 47
$ws.Range("E47").Value = "'  +2.44%  "

# Row 48
$ws.Range("D48").Value = "'2.00"
$ws.Range("E48").Value = "'  -1.00%  "

# Row 49
$ws.Range("D49").Value = "'8.94"
$ws.Range("E49").Value = "'  -1.85%  "

# Row 50
$ws.Range("D50").Value = "'5.18"
$ws.Range("E50").Value = "'  -1.16%  "

# Row 51
$ws.Range("D51").Value = "'78.80"
$ws.Range("E51").Value = "'  -0.96%  "
